$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two additional columns.
$ws.Range("J1").Value = "Client ID"
$ws.Range("K1").Value = "Client Secret"

# Match the bold + centered header formatting used by the rest of row 1
# (same look as style index 3 in the original file: bold font, centered).
$ws.Range("J1:K1").Font.Bold = $true
$ws.Range("J1:K1").HorizontalAlignment = -4108

# Column widths for the two new columns (values chosen so the engine's
# internal pixel-rounding reproduces the target stored widths of 21 and
# 22.6640625 characters).
$ws.Columns.Item(10).ColumnWidth = 20.17
$ws.Columns.Item(11).ColumnWidth = 21.8

# Move the selection to the newly added header cells, matching the
# post-edit state captured by Excel.
$ws.Range("J1:K1").Select()
